$d = $word.ActiveDocument

$pairs = @(
    ,@('72-60=', '37-16=')
    ,@('5+81=', '99-48=')
    ,@('38+15=', '8+66=')
    ,@('88-24=', '86-63=')
    ,@('44+20=', '17+81=')
    ,@('72-3=', '8+85=')
    ,@('24-11=', '49+10=')
    ,@('89-42=', '38+35=')
    ,@('30+19=', '3+64=')
    ,@('28+8=', '48-1=')
    ,@('16+60=', '51+35=')
    ,@('98-49=', '77+13=')
    ,@('57+20=', '6+71=')
    ,@('9+17=', '6+73=')
    ,@('84-80=', '14+61=')
    ,@('53-29=', '98-50=')
    ,@('10+20=', '45+17=')
    ,@('5+54=', '1+5=')
    ,@('65-61=', '38+13=')
    ,@('96-86=', '22+73=')
    ,@('80-3=', '59+14=')
    ,@('33-27=', '41+49=')
    ,@('79-16=', '30+59=')
    ,@('41+13=', '64-59=')
    ,@('17+17=', '56-10=')
    ,@('91-89=', '56+20=')
    ,@('71-61=', '6+12=')
    ,@('1+60=', '42-24=')
    ,@('0+66=', '90-72=')
    ,@('41+2=', '81-25=')
    ,@('81-63=', '67-5=')
    ,@('55+17=', '2+12=')
    ,@('2+42=', '66+21=')
    ,@('39+45=', '73+22=')
    ,@('22-9=', '86-36=')
    ,@('88-30=', '86-7=')
    ,@('42+44=', '17+36=')
    ,@('80-35=', '76-61=')
    ,@('44+28=', '57-53=')
    ,@('72+8=', '9+6=')
    ,@('52+0=', '32-19=')
    ,@('1+85=', '8+51=')
    ,@('43+48=', '49-35=')
    ,@('47+0=', '94-88=')
    ,@('47+27=', '72+18=')
    ,@('58-7=', '86-26=')
    ,@('38-15=', '74-4=')
    ,@('63+5=', '3+63=')
    ,@('4+12=', '79-58=')
    ,@('3+4=', '49+6=')
    ,@('91-7=', '60-53=')
    ,@('49-17=', '6+7=')
    ,@('68-66=', '47-28=')
    ,@('50-0=', '32-16=')
    ,@('63-11=', '43-19=')
    ,@('37-10=', '28+48=')
    ,@('44+13=', '92-78=')
    ,@('39+32=', '67-43=')
    ,@('89-24=', '42+30=')
    ,@('24+40=', '95-16=')
    ,@('25+40=', '59+9=')
    ,@('18+75=', '97-44=')
    ,@('27-18=', '78-55=')
    ,@('62-41=', '17+47=')
    ,@('30+13=', '51-29=')
    ,@('91-9=', '22+71=')
    ,@('2+30=', '57-24=')
    ,@('6+88=', '17-0=')
    ,@('3+53=', '14+72=')
    ,@('57-36=', '46-26=')
    ,@('5+10=', '71-6=')
    ,@('55+35=', '40+37=')
    ,@('50+10=', '86-52=')
    ,@('8+69=', '22+18=')
    ,@('52-20=', '42-16=')
    ,@('67-22=', '38+3=')
    ,@('65-26=', '63-19=')
    ,@('31+21=', '24+52=')
    ,@('4+45=', '6+17=')
    ,@('80-17=', '86-55=')
    ,@('53+2=', '83+2=')
    ,@('28+31=', '34+2=')
    ,@('5+47=', '91-53=')
    ,@('71+24=', '45+9=')
    ,@('71-57=', '99-72=')
    ,@('10+27=', '57+14=')
    ,@('1+45=', '69-17=')
    ,@('52-23=', '15+18=')
    ,@('79-17=', '72-67=')
    ,@('83-55=', '99-68=')
    ,@('26-18=', '10+67=')
    ,@('31+14=', '1+49=')
    ,@('12+40=', '89-64=')
    ,@('76-22=', '73-72=')
    ,@('67-0=', '88+3=')
    ,@('8+72=', '5+39=')
    ,@('83-71=', '44+8=')
    ,@('16-12=', '90-20=')
    ,@('39+20=', '28+52=')
    ,@('92-60=', '5+83=')
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done."
